$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the CodeSystem version (0.4.0 -> 0.7.0)
$ws.Range("B3").Value = "0.7.0"

# Remove the Jurisdiction / Chile row entirely; rows below shift up
$ws.Rows.Item(11).Delete()
